# Ticket 79 - Fix implicit sheet cloning case when the number of items in
# the collection is 1.
#
# The "ImplCloningNormalTemplate.xlsx" template demonstrates JETT's implicit
# sheet cloning feature: the sheet named "${dvs.name}$@l=0" is a template
# sheet that gets cloned once per item in a collection whose size is used as
# the "limit" (l=...) in the sheet-name directive. To exercise (and regress
# test) the case where the collection has more than one item, a second demo
# sheet "${dvs.name}$@l=1" is added as an exact copy of the original
# "${dvs.name}$@l=0" sheet, placed right after it (at the end of the
# workbook).

$wb = $excel.ActiveWorkbook

# The template sheet we are cloning. NOTE: must be single-quoted -- this
# sheet name contains "${...}" sequences that PowerShell would otherwise try
# to expand as variables inside a double-quoted string.
$sourceName = '${dvs.name}$@l=0'
$newName    = '${dvs.name}$@l=1'

$src = $wb.Worksheets.Item($sourceName)

# Copy the source sheet to the very end of the workbook (after the last
# existing sheet, i.e. "Static3").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy([System.Reflection.Missing]::Value, $lastSheet)

# The freshly copied sheet is now the last sheet in the workbook; give it
# its proper template name.
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = $newName

# Restore/normalize the selections: the source sheet ends up with A1:E4
# selected (its full used range), while the new clone keeps the original
# A1:E1 selection (the merged title cell) that the source sheet had before
# being copied.
$src.Range("A1:E4").Select() | Out-Null

# Leave the workbook's original active sheet ("Static1") selected/active,
# same as before the edit.
$wb.Worksheets.Item("Static1").Activate() | Out-Null
